$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TRAINING LOG")
$ws2 = $wb.Worksheets.Item("PREDICTION LOG")

# --- Sheet1: TRAINING LOG ---
$trainRows = @(
    @(4, "BINARY", "24:48:1024", "2021-09-06-23-03-04", 1024, "ABC"),
    @(5, "BINARY", "24:48:1024", "2021-09-06-23-03-09", 1024, "ABC"),
    @(6, "BINARY", "24:48:1024", "2021-09-06-23-04-34", 1024, "ABC")
)

$r = 7
foreach ($row in $trainRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# --- Sheet2: PREDICTION LOG ---
$predRows = @(
    @(4, "BINARY", "24:48:1024", "2021-09-06-23-03-04", 1024, "ABC", "2021-09-06-23-03-04", "DEF.txt", 0.96, 0.9099),
    @(5, "BINARY", "24:48:1024", "2021-09-06-23-03-09", 1024, "ABC", "2021-09-06-23-03-09", "DEF.txt", 0.96, 0.9099),
    @(6, "BINARY", "24:48:1024", "2021-09-06-23-04-34", 1024, "ABC", "2021-09-06-23-04-34", "DEF.txt", 0.96, 0.9099)
)

$r = 7
foreach ($row in $predRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 10).Value = $row[6]
    $ws2.Cells.Item($r, 11).Value = $row[7]
    $ws2.Cells.Item($r, 12).Value = $row[8]
    $ws2.Cells.Item($r, 13).Value = $row[9]
    $r++
}
